$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.115.29"
$ws.Range("E2").Value = "  -2.12%  "
Set-TextValue $ws.Range("D3") "1.856.47"
$ws.Range("E3").Value = "  -3.56%  "
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  +0.10%  "
Set-TextValue $ws.Range("D5") "233.73"
$ws.Range("E5").Value = "  -3.17%  "
Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  +0.17%  "
Set-TextValue $ws.Range("D7") "0.4653"
$ws.Range("E7").Value = "  -2.73%  "
Set-TextValue $ws.Range("D8") "0.2817"
$ws.Range("E8").Value = "  -2.25%  "
Set-TextValue $ws.Range("D9") "0.06553"
$ws.Range("E9").Value = "  -3.33%  "
Set-TextValue $ws.Range("D10") "19.95"
$ws.Range("E10").Value = "  +1.80%  "
Set-TextValue $ws.Range("D11") "0.07821"
$ws.Range("E11").Value = "  +0.49%  "
Set-TextValue $ws.Range("D12") "96.70"
$ws.Range("E12").Value = "  -7.32%  "
Set-TextValue $ws.Range("D13") "1.866.25"
$ws.Range("E13").Value = "  -3.12%  "
Set-TextValue $ws.Range("D14") "5.105"
$ws.Range("E14").Value = "  -3.22%  "
Set-TextValue $ws.Range("D15") "0.6643"
$ws.Range("E15").Value = "  -2.53%  "
Set-TextValue $ws.Range("D16") "282.04"
$ws.Range("E16").Value = "  -3.02%  "
Set-TextValue $ws.Range("D17") "30.156.96"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  -1.09%  "
Set-TextValue $ws.Range("D20") "12.59"
$ws.Range("E20").Value = "  -2.14%  "
Set-TextValue $ws.Range("D21") "2.108.69"
$ws.Range("E21").Value = "  -3.18%  "
Set-TextValue $ws.Range("D22") "0.000007225"
$ws.Range("E22").Value = "  -4.65%  "
$ws.Range("E23").Value = "  +0.02%  "
Set-TextValue $ws.Range("D24") "6.131"
$ws.Range("E24").Value = "  -3.89%  "
Set-TextValue $ws.Range("D25") "167.88"
$ws.Range("E25").Value = "  -0.03%  "
Set-TextValue $ws.Range("D26") "9.319"
$ws.Range("E26").Value = "  -2.13%  "
Set-TextValue $ws.Range("D27") "18.86"
$ws.Range("E27").Value = "  -4.42%  "
Set-TextValue $ws.Range("D28") "1.911"
$ws.Range("E28").Value = "  -9.57%  "
Set-TextValue $ws.Range("D29") "1.334"
$ws.Range("E29").Value = "  -3.98%  "
Set-TextValue $ws.Range("D30") "0.09561"
$ws.Range("E30").Value = "  -5.01%  "
Set-TextValue $ws.Range("D32") "1.470"
$ws.Range("E32").Value = "  -3.75%  "
Set-TextValue $ws.Range("D33") "4.099"
$ws.Range("E33").Value = "  -5.10%  "
Set-TextValue $ws.Range("D34") "0.04650"
$ws.Range("E34").Value = "  -3.30%  "
Set-TextValue $ws.Range("D35") "0.7005"
$ws.Range("E35").Value = "  -4.72%  "
Set-TextValue $ws.Range("D36") "1.097"
$ws.Range("E36").Value = "  -2.44%  "
Set-TextValue $ws.Range("D37") "2.702"
$ws.Range("E37").Value = "  -0.44%  "
Set-TextValue $ws.Range("D38") "0.01846"
$ws.Range("E38").Value = "  -4.90%  "
$ws.Range("E39").Value = "  -1.24%  "
Set-TextValue $ws.Range("D40") "2.511"
$ws.Range("E40").Value = "  -4.14%  "
Set-TextValue $ws.Range("D41") "72.07"
$ws.Range("E41").Value = "  -3.90%  "
Set-TextValue $ws.Range("D42") "0.8519"
$ws.Range("E42").Value = "  -1.81%  "
Set-TextValue $ws.Range("D43") "1.924"
$ws.Range("E43").Value = "  -4.66%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D45") "103.84"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D46") "0.4142"
$ws.Range("E46").Value = "  -4.33%  "
Set-TextValue $ws.Range("D47") "989.60"
$ws.Range("E47").Value = "  +0.74%  "
Set-TextValue $ws.Range("D48") "7.223"
$ws.Range("E48").Value = "  -4.35%  "
Set-TextValue $ws.Range("D49") "9.201"
$ws.Range("E49").Value = "  +2.81%  "
Set-TextValue $ws.Range("D50") "33.93"
$ws.Range("E50").Value = "  -2.81%  "
Set-TextValue $ws.Range("D51") "0.1139"
$ws.Range("E51").Value = "  -5.92%  "
